$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells in this sheet store their numbers as literal text (t="inlineStr"),
# not real numeric values, so a leading apostrophe is used to force
# Excel's text interpretation instead of silently converting the
# numeric-looking strings into Number cells.
$ws.Range("C10").Value = "'239.89"
$ws.Range("D10").Value = "'403.89"

$ws.Range("C16").Value = "'2,541.95"
$ws.Range("D16").Value = "'2,982.39"

$ws.Range("C17").Value = "'886.41"
$ws.Range("D17").Value = "'1,383.31"

$ws.Range("C24").Value = "'886.41"
$ws.Range("D24").Value = "'1,383.31"
